$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "as of" timestamp (column D) for all data rows 2-56
$ws.Range("D2:D56").Value = 45992.288842592592

# Refresh rows 19-56 with the latest report data (A: station, B: terminal, C: last charge end time)
$ws.Range("A19").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B19").Value = '502号直流'
$ws.Range("C19").Value = 45989.209733796299
$ws.Range("A20").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B20").Value = '603号直流'
$ws.Range("C20").Value = 45989.545405092591
$ws.Range("A21").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B21").Value = '805号直流'
$ws.Range("C21").Value = 45989.926840277774
$ws.Range("A22").Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Range("B22").Value = '406号直流'
$ws.Range("C22").Value = 45990.197141203702
$ws.Range("A23").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B23").Value = '905号直流'
$ws.Range("C23").Value = 45990.737581018519
$ws.Range("A24").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B24").Value = '505号直流'
$ws.Range("C24").Value = 45991.02715277778
$ws.Range("A25").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B25").Value = 'B02号直流'
$ws.Range("C25").Value = 45991.033067129632
$ws.Range("A26").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B26").Value = '702号直流'
$ws.Range("C26").Value = 45991.033333333333
$ws.Range("A27").Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Range("B27").Value = '103号直流'
$ws.Range("C27").Value = 45991.050717592596
$ws.Range("A28").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B28").Value = '904号直流'
$ws.Range("C28").Value = 45991.157314814816
$ws.Range("A29").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B29").Value = '902号直流'
$ws.Range("C29").Value = 45991.186979166669
$ws.Range("A30").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B30").Value = '703号直流'
$ws.Range("C30").Value = 45991.200173611112
$ws.Range("A31").Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Range("B31").Value = '103号直流'
$ws.Range("C31").Value = 45991.246099537035
$ws.Range("A32").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B32").Value = '402号直流'
$ws.Range("C32").Value = 45991.264398148145
$ws.Range("A33").Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Range("B33").Value = '105号直流'
$ws.Range("C33").Value = 45991.451932870368
$ws.Range("A34").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B34").Value = '301号直流'
$ws.Range("C34").Value = 45991.493321759262
$ws.Range("A35").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B35").Value = '203号直流'
$ws.Range("C35").Value = 45991.515543981484
$ws.Range("A36").Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Range("B36").Value = '203号直流'
$ws.Range("C36").Value = 45991.539375
$ws.Range("A37").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B37").Value = '105号直流'
$ws.Range("C37").Value = 45991.544074074074
$ws.Range("A38").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B38").Value = '111号直流'
$ws.Range("C38").Value = 45991.552152777775
$ws.Range("A39").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B39").Value = '207号直流'
$ws.Range("C39").Value = 45991.553101851852
$ws.Range("A40").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B40").Value = '110号直流'
$ws.Range("C40").Value = 45991.559398148151
$ws.Range("A41").Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Range("B41").Value = '003B号直流'
$ws.Range("C41").Value = 45991.563310185185
$ws.Range("A42").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B42").Value = '705号直流'
$ws.Range("C42").Value = 45991.58258101852
$ws.Range("A43").Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Range("B43").Value = '905号直流'
$ws.Range("C43").Value = 45991.594629629632
$ws.Range("A44").Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Range("B44").Value = '102号直流'
$ws.Range("C44").Value = 45991.600381944445
$ws.Range("A45").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B45").Value = 'A01号直流'
$ws.Range("C45").Value = 45991.616446759261
$ws.Range("A46").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B46").Value = '311号直流'
$ws.Range("C46").Value = 45991.635277777779
$ws.Range("A47").Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Range("B47").Value = '903号直流'
$ws.Range("C47").Value = 45991.644155092596
$ws.Range("A48").Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Range("B48").Value = '001B号直流'
$ws.Range("C48").Value = 45991.653356481482
$ws.Range("A49").Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Range("B49").Value = '008B号直流'
$ws.Range("C49").Value = 45991.657337962963
$ws.Range("A50").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B50").Value = '903号直流'
$ws.Range("C50").Value = 45991.673692129632
$ws.Range("A51").Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Range("B51").Value = '306号直流'
$ws.Range("C51").Value = 45991.676076388889
$ws.Range("A52").Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Range("B52").Value = '404号直流'
$ws.Range("C52").Value = 45991.676944444444
$ws.Range("A53").Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Range("B53").Value = 'A02号直流'
$ws.Range("C53").Value = 45991.746689814812
$ws.Range("A54").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B54").Value = '108号直流'
$ws.Range("C54").Value = 45991.754513888889
$ws.Range("A55").Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Range("B55").Value = '305号直流'
$ws.Range("C55").Value = 45991.774409722224
$ws.Range("A56").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B56").Value = '206号直流'
$ws.Range("C56").Value = 45991.778541666667

# Restore the active-cell selection shown when the workbook was last saved
$ws.Range("E10").Select()
